$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row (62) with the latest "cotações" (quotes) for 2025-11-05,
# following the same layout as the existing rows (A = date serial, B:E = text values)

$newRow = 62

# Copy the date cell's formatting from the previous row so the new date cell
# matches the existing date column formatting
$ws.Range("A61").Copy()
$ws.Range("A62").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("A" + $newRow).Value = 45966
$ws.Range("B" + $newRow).Value = "22,0559"
$ws.Range("C" + $newRow).Value = "11,2354"
$ws.Range("D" + $newRow).Value = "15,3884"
$ws.Range("E" + $newRow).Value = "15,3884"

$wb.Save()
